$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B3 text: "240 + 143" -> "(240 + 143)"
$ws.Range("B3").Value = "(240 + 143)"

# Apply a numeric format (2 decimal places) to B2:B4
$ws.Range("B2:B4").NumberFormat = "0.00"

# Set column B width (closest achievable match for the target stored width of 18.5703125)
$ws.Columns("B").ColumnWidth = 17.7

# Update selection to C8
$ws.Range("C8").Select()
